# TC_65166 - Updated test data for DC, TripCurrent, Voltdrop, BatteryStandby
# Apply the same shape of edit to both sheets: add a "Loop"/"Column" mini
# table in columns H:I, tweak a few DC-unit numbers, and stamp a
# requirement id into B4.

$wb = $excel.ActiveWorkbook

foreach ($idx in 1,2) {
    $ws = $wb.Worksheets.Item($idx)

    # --- Copy existing formatting onto the new H:I "Loop"/"Column" cells
    # before we touch their source cells' own content/format.
    $ws.Range("E1").Copy()
    $ws.Range("H1:I1").PasteSpecial(-4122)   # xlPasteFormats

    $ws.Range("B4").Copy()
    $ws.Range("H2:H5").PasteSpecial(-4122)   # xlPasteFormats
    $ws.Range("I2").PasteSpecial(-4122)      # xlPasteFormats

    # --- New "Loop"/"Column" reference table in columns H & I
    $ws.Range("H1").Value = "Loop"
    $ws.Range("I1").Value = "Column"
    $ws.Range("H2").Value = "Built-in Loop-A"
    $ws.Range("I2").Value = 2
    $ws.Range("H3").Value = "Built-in Loop-B"
    $ws.Range("H4").Value = "Built-in Loop-C"
    $ws.Range("H5").Value = "Built-in Loop-D"

    # --- Stamp the requirement / defect id into B4 (formatting cleared)
    $ws.Range("B4").ClearFormats()
    $ws.Range("B4").Value = "NGC-1826/T920 OR TC-65166"
}

# --- Sheet-specific DC Unit number tweaks
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("G1").Value = 300
$ws1.Range("G2").Value = "'"      # re-enter as blank/quote-prefixed cell
$ws1.Range("G5").Value = 340

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("G1").Value = 340
$ws2.Range("G2").Value = 319
$ws2.Range("G3").Value = 343
$ws2.Range("G5").Value = 360

# --- Restore selection: both sheets now have B4 selected, sheet1 stays
# the active tab.
$ws2.Activate() | Out-Null
$ws2.Range("B4").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("B4").Select() | Out-Null
